# Generate Report for Handoff
# Swap the old run id (69c5e5e8-904f-498e-81e3-9491c1dc6c9c) for the new one
# (62c30c64-a7b2-41a5-b56a-f7ee2f0bdfcb), refresh the xliff checksum suffix,
# and bump the handoff timestamps on all three sheets.

$wb = $excel.ActiveWorkbook

$oldId = "69c5e5e8-904f-498e-81e3-9491c1dc6c9c"
$newId = "62c30c64-a7b2-41a5-b56a-f7ee2f0bdfcb"
$oldHash = "333e4cecc2f175c2a068144a5f71a3cd93e0aa52"
$newHash = "968898bc43c27254670ca6d8c887bfef0dd0bb6b"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d59e3e7a12c5ee2431dc8557a938549cb570b402/e2e/$oldId.md"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"

$pathAndName = "e2e\$newId.md"
$rngB2 = $wsOverview.Range("B2")
$rngB2.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($rngB2, $hyperlinkAddress, $null, $null, $pathAndName)
$rngB2.Value = $pathAndName

$wsOverview.Range("G2").Value = "2016-08-18 04:57:52"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$rngZhA2 = $wsZh.Range("A2")
$rngZhA2.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($rngZhA2, $hyperlinkAddress, $null, $null, "$newId.md")
$rngZhA2.Value = "$newId.md"

$wsZh.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-18 04:57:47"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$rngDeA2 = $wsDe.Range("A2")
$rngDeA2.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($rngDeA2, $hyperlinkAddress, $null, $null, "$newId.md")
$rngDeA2.Value = "$newId.md"

$wsDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-18 04:57:52"
